$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.24291338728751896
$ws.Range("C2").Value = 0.1577086122646325
$ws.Range("D2").Value = 0.3281181623104054
$ws.Range("E2").Value = -0.01315163965455784
$ws.Range("F2").Value = -0.024757441220282374
$ws.Range("G2").Value = -0.0034830979512390876
$ws.Range("K2").Value = 0.16341024123474443
$ws.Range("L2").Value = 0.07140362280617438
$ws.Range("M2").Value = 0.2554168596633145
$ws.Range("K3").Value = -0.27848725286618436
$ws.Range("L3").Value = -0.3571216353923619
$ws.Range("M3").Value = -0.19985287034000682
$ws.Range("B4").Value = 0.21485099728507143
$ws.Range("C4").Value = 0.14923864702015535
$ws.Range("D4").Value = 0.2804633475499875
$ws.Range("E4").Value = -0.010264346445229949
$ws.Range("F4").Value = -0.01863872606367311
$ws.Range("G4").Value = -0.002686931767029087
$ws.Range("K4").Value = 0.10262579001870845
$ws.Range("L4").Value = 0.027223560842243172
$ws.Range("M4").Value = 0.17802801919517372
$ws.Range("B5").Value = 0.18469761295441794
$ws.Range("C5").Value = 0.09784387807518917
$ws.Range("D5").Value = 0.2715513478336467
$ws.Range("E5").Value = -0.012116228849009099
$ws.Range("F5").Value = -0.02261648439768467
$ws.Range("G5").Value = -0.0032881203977300176
$ws.Range("K5").Value = 0.15296185129777135
$ws.Range("L5").Value = 0.06721605338646187
$ws.Range("M5").Value = 0.23870764920908083
$ws.Range("B6").Value = 0.13438256788187808
$ws.Range("C6").Value = 0.0637781043590879
$ws.Range("D6").Value = 0.20498703140466826
$ws.Range("E6").Value = -0.009319932775022344
$ws.Range("F6").Value = -0.01700283746050362
$ws.Range("G6").Value = -0.002430282245406993
$ws.Range("B7").Value = 0.15591409345080884
$ws.Range("C7").Value = 0.09372872545921257
$ws.Range("D7").Value = 0.2180994614424051
$ws.Range("E7").Value = -0.010129456808685235
$ws.Range("F7").Value = -0.017829874641662406
$ws.Range("G7").Value = -0.003181235916527662
$ws.Range("H7").Value = 0.0058142459436786005
$ws.Range("I7").Value = 0.0012877890824605732
$ws.Range("J7").Value = 0.010043374554235636
$ws.Range("K7").Value = 0.1925913397006869
$ws.Range("L7").Value = 0.12818029232186975
$ws.Range("M7").Value = 0.25700238707950407
$ws.Range("B8").Value = 0.10694518928460074
$ws.Range("C8").Value = 0.051301157275299575
$ws.Range("D8").Value = 0.1625892212939019
$ws.Range("E8").Value = -0.00763166263540419
$ws.Range("F8").Value = -0.014308957107606175
$ws.Range("G8").Value = -0.0021400137690131054
$ws.Range("H8").Value = 0.0010725906128833957
$ws.Range("I8").Value = [double]"-4.411110440798947e-05"
$ws.Range("J8").Value = 0.0032913802373276957
$ws.Range("K8").Value = 0.1911970231649376
$ws.Range("L8").Value = 0.1162203456560478
$ws.Range("M8").Value = 0.2661737006738274
